# Applies the edits described by the diff:
#  - Adds a "Medida" header label in B4
#  - Re-enters the measured values (columns C, D, E) and the computed
#    statistics (columns L, M, N) for rows 5-11 using a comma typed as an
#    apostrophe ( ' ) instead of a decimal point, which makes Excel store
#    them as text instead of numbers (matching the target workbook).
#  - Moves the active selection to D12 (was I13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the "Medida" column
$ws.Range("B4").Value = "Medida"

# Row 5
$ws.Range("D5").Value = "14'6"
$ws.Range("M5").Value = "0'489"

# Row 6
$ws.Range("C6").Value = "2'6"
$ws.Range("D6").Value = "23'8"
$ws.Range("E6").Value = "19'4"
$ws.Range("L6").Value = "1'516"
$ws.Range("M6").Value = "0'748"
$ws.Range("N6").Value = "0'489"

# Row 7
$ws.Range("C7").Value = "3'2"
$ws.Range("D7").Value = "30'8"
$ws.Range("E7").Value = "32'4"
$ws.Range("L7").Value = "0'447"
$ws.Range("M7").Value = "0'748"
$ws.Range("N7").Value = "0'8"

# Row 8
$ws.Range("C8").Value = "1'2"
$ws.Range("D8").Value = "41'8"
$ws.Range("E8").Value = "41'2"
$ws.Range("L8").Value = "1'095"
$ws.Range("M8").Value = "0'4"
$ws.Range("N8").Value = "1'469"

# Row 9
$ws.Range("C9").Value = "1'2"
$ws.Range("D9").Value = "49'6"
$ws.Range("E9").Value = "50'6"
$ws.Range("L9").Value = "0'836"
$ws.Range("M9").Value = "0'489"
$ws.Range("N9").Value = "8'522"

# Row 10
$ws.Range("C10").Value = "1'2"
$ws.Range("D10").Value = "62'8"
$ws.Range("E10").Value = "60'4"
$ws.Range("L10").Value = "0'447"
$ws.Range("M10").Value = "0'4"
$ws.Range("N10").Value = "2'244"

# Row 11
$ws.Range("C11").Value = "0'8"
$ws.Range("D11").Value = "71'0"
$ws.Range("E11").Value = "65'6"
$ws.Range("L11").Value = "0'836"
$ws.Range("N11").Value = "2'939"

# Move the selection, matching the saved cursor position in the edited file
$ws.Range("D12").Select()
